# New Test Case Added.
# Adds a new "Invalid login" test case row (row 6) and fills in the
# previously-empty row 7 ("demo"/"test") with an Invalid result.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6: new test case - username "testfail" / password "failtest" -> Invalid
$ws.Range("A6").Value = "testfail"
$ws.Range("B6").Value = "failtest"
$ws.Range("C6").Value = "Invalid"

# Row 7: new test case - username "demo" / password "test" -> Invalid
$ws.Range("A7").Value = "demo"
$ws.Range("B7").Value = "test"
$ws.Range("C7").Value = "Invalid"

# C7 already carries the bordered "result" cell formatting; copy that
# formatting onto the newly written C6 so both result cells match.
$ws.Range("C7").Copy()
$ws.Range("C6").PasteSpecial(-4122)  # xlPasteFormats

# Leave the selection on C7, matching the saved worksheet view.
$ws.Range("C7").Select()
